$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Conference" header to "Conf."
$ws.Range("B1").Value = "Conf."

# Swap the order of "NET" and "Record" columns so that
# "Record" comes before "NET" (C1 = Record, D1 = NET)
$ws.Range("C1").Value = "Record"
$ws.Range("D1").Value = "NET"
